# Update DM integration test fixture
# 1. Change the CodeSchemes ID value (A2) on sheet "CodeSchemes"
# 2. Bold the header row on all three sheets
# 3. Widen the columns (reflecting the width Excel computes once the header
#    text becomes bold) on all three sheets

$wb = $excel.ActiveWorkbook

$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodes       = $wb.Worksheets.Item("Codes")
$wsExtensions  = $wb.Worksheets.Item("Extensions")

# --- 1. Update the code scheme identifier -------------------------------
$wsCodeSchemes.Range("A2").Value = "1824d625-7eec-490e-885e-e00b903f28ac"

# --- 2. Bold the header rows ----------------------------------------------
# Each sheet gets its own (otherwise identical) bold style, matching the
# three separate font/xf entries the source workbook ends up with.
$boldStyle1 = $wb.Styles.Add("HeaderBold_CodeSchemes")
$boldStyle1.Font.Bold = $true
$wsCodeSchemes.Range("A1:N1").Style = "HeaderBold_CodeSchemes"

$boldStyle2 = $wb.Styles.Add("HeaderBold_Codes")
$boldStyle2.Font.Bold = $true
$wsCodes.Range("A1:J1").Style = "HeaderBold_Codes"

$boldStyle3 = $wb.Styles.Add("HeaderBold_Extensions")
$boldStyle3.Font.Bold = $true
$wsExtensions.Range("A1:I1").Style = "HeaderBold_Extensions"

# --- 3. Resize the columns to match the new (bold) best-fit widths -------

# CodeSchemes
$wsCodeSchemes.Columns.Item(1).ColumnWidth  = 33.428571428571429
$wsCodeSchemes.Columns.Item(2).ColumnWidth  = 17.428571428571427
$wsCodeSchemes.Columns.Item(3).ColumnWidth  = 25.714285714285715
$wsCodeSchemes.Columns.Item(4).ColumnWidth  = 22.428571428571427
$wsCodeSchemes.Columns.Item(5).ColumnWidth  = 14.142857142857142
$wsCodeSchemes.Columns.Item(6).ColumnWidth  = 19.142857142857142
$wsCodeSchemes.Columns.Item(7).ColumnWidth  = 20.142857142857142
$wsCodeSchemes.Columns.Item(8).ColumnWidth  = 19.142857142857142
$wsCodeSchemes.Columns.Item(9).ColumnWidth  = 20.714285714285715
$wsCodeSchemes.Columns.Item(10).ColumnWidth = 24.0
$wsCodeSchemes.Columns.Item(11).ColumnWidth = 19.142857142857142
$wsCodeSchemes.Columns.Item(12).ColumnWidth = 15.714285714285714
$wsCodeSchemes.Columns.Item(13).ColumnWidth = 20.714285714285715
$wsCodeSchemes.Columns.Item(14).ColumnWidth = 27.285714285714285

# Codes
$wsCodes.Columns.Item(1).ColumnWidth  = 5.857142857142857
$wsCodes.Columns.Item(2).ColumnWidth  = 17.428571428571427
$wsCodes.Columns.Item(3).ColumnWidth  = 15.714285714285714
$wsCodes.Columns.Item(4).ColumnWidth  = 14.142857142857142
$wsCodes.Columns.Item(5).ColumnWidth  = 15.714285714285714
$wsCodes.Columns.Item(6).ColumnWidth  = 19.142857142857142
$wsCodes.Columns.Item(7).ColumnWidth  = 20.714285714285715
$wsCodes.Columns.Item(8).ColumnWidth  = 24.0
$wsCodes.Columns.Item(9).ColumnWidth  = 19.142857142857142
$wsCodes.Columns.Item(10).ColumnWidth = 15.714285714285714

# Extensions
$wsExtensions.Columns.Item(1).ColumnWidth = 5.857142857142857
$wsExtensions.Columns.Item(2).ColumnWidth = 17.428571428571427
$wsExtensions.Columns.Item(3).ColumnWidth = 14.142857142857142
$wsExtensions.Columns.Item(4).ColumnWidth = 24.0
$wsExtensions.Columns.Item(5).ColumnWidth = 15.714285714285714
$wsExtensions.Columns.Item(6).ColumnWidth = 19.142857142857142
$wsExtensions.Columns.Item(7).ColumnWidth = 19.142857142857142
$wsExtensions.Columns.Item(8).ColumnWidth = 15.714285714285714
$wsExtensions.Columns.Item(9).ColumnWidth = 24.0
